$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.030389785766602
$ws.Range("B1").Value = 1.836334943771362
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.772303819656372
$ws.Range("E1").Value = 1.231554388999939
